$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"
$ws.Range("B7").Value = "Pabellón De Arteaga"
$ws.Range("B8").Value = "Rincón De Romos"
$ws.Range("B9").Value = "San José De Gracia"
$ws.Range("B33").Value = "Comitán De Domínguez"
$ws.Range("B48").Value = "San Cristóbal De Las Casas"
$ws.Range("B71").Value = "Coyame Del Sotol"
$ws.Range("B81").Value = "Guadalupe Y Calvo"
$ws.Range("B84").Value = "Hidalgo Del Parral"
$ws.Range("B106").Value = "San Francisco De Borja"
$ws.Range("B107").Value = "San Francisco De Conchos"
$ws.Range("B108").Value = "San Francisco Del Oro"
$ws.Range("B116").Value = "Valle De Zaragoza"
$ws.Range("A118").Value = "Ciudad De México"
$ws.Range("A134").Value = "Coahuila De Zaragoza"
$ws.Range("B154").Value = "Villa De Álvarez"
$ws.Range("B158").Value = "Coneto De Comonfort"
$ws.Range("B172").Value = "Nombre De Dios"
$ws.Range("B178").Value = "Pánuco De Coronado"
$ws.Range("B182").Value = "San Juan De Guadalupe"
$ws.Range("B183").Value = "San Juan Del Río"
$ws.Range("B184").Value = "San Luis Del Cordero"
$ws.Range("A194").Value = "Estado De México"
$ws.Range("B195").Value = "Atizapán De Zaragoza"
$ws.Range("B201").Value = "Ecatepec De Morelos"
$ws.Range("B208").Value = "Naucalpan De Juárez"
$ws.Range("B216").Value = "Tenango Del Valle"
$ws.Range("B217").Value = "Tlalnepantla De Baz"
$ws.Range("B219").Value = "Valle De Chalco Solidaridad"
$ws.Range("B228").Value = "Dolores Hidalgo Cuna De La Independencia Nacional"
$ws.Range("B232").Value = "Jaral Del Progreso"
$ws.Range("B239").Value = "San Diego De La Unión"
$ws.Range("B241").Value = "San Francisco Del Rincón"
$ws.Range("B242").Value = "San Luis De La Paz"
$ws.Range("B243").Value = "San Miguel De Allende"
$ws.Range("B244").Value = "Silao De La Victoria"
$ws.Range("B249").Value = "Valle De Santiago"
$ws.Range("B252").Value = "Acapulco De Juárez"
$ws.Range("B253").Value = "Ajuchitlán Del Progreso"
$ws.Range("B255").Value = "Atoyac De Álvarez"
$ws.Range("B256").Value = "Ayutla De Los Libres"
$ws.Range("B258").Value = "Chilapa De Álvarez"
$ws.Range("B259").Value = "Chilpancingo De Los Bravo"
$ws.Range("B260").Value = "Coyuca De Catalán"
$ws.Range("B261").Value = "Cuetzala Del Progreso"
$ws.Range("B262").Value = "Cutzamala De Pinzón"
$ws.Range("B265").Value = "Huitzuco De Los Figueroa"
$ws.Range("B266").Value = "Iguala De La Independencia"
$ws.Range("B268").Value = "La Unión De Isidoro Montes De Oca"
$ws.Range("B272").Value = "Taxco De Alarcón"
$ws.Range("B274").Value = "Tepecoacuilco De Trujano"
$ws.Range("B275").Value = "Técpan De Galeana"
$ws.Range("B277").Value = "Zihuatanejo De Azueta"
$ws.Range("B287").Value = "Progreso De Obregón"
$ws.Range("B290").Value = "Tula De Allende"
$ws.Range("B291").Value = "Tulancingo De Bravo"
$ws.Range("B294").Value = "Ahualulco De Mercado"
$ws.Range("B298").Value = "Atotonilco El Alto"
$ws.Range("B299").Value = "Autlán De Navarro"
$ws.Range("B307").Value = "Encarnación De Díaz"
$ws.Range("B309").Value = "Huejuquilla El Alto"
$ws.Range("B311").Value = "Ixtlahuacán Del Río"
$ws.Range("B314").Value = "Lagos De Moreno"
$ws.Range("B321").Value = "San Diego De Alejandría"
$ws.Range("B323").Value = "San Juan De Los Lagos"
$ws.Range("B325").Value = "San Miguel El Alto"
$ws.Range("B328").Value = "Talpa De Allende"
$ws.Range("B329").Value = "Tamazula De Gordiano"
$ws.Range("B333").Value = "Teocuitatlán De Corona"
$ws.Range("B334").Value = "Tepatitlán De Morelos"
$ws.Range("B335").Value = "Tizapán El Alto"
$ws.Range("B336").Value = "Tlajomulco De Zúñiga"
$ws.Range("B338").Value = "Unión De Tula"
$ws.Range("B341").Value = "Yahualica De González Gallo"
$ws.Range("B344").Value = "Zapotlán Del Rey"
$ws.Range("B345").Value = "Zapotlán El Grande"
$ws.Range("A347").Value = "Michoacán De Ocampo"
$ws.Range("B395").Value = "Tetela Del Volcán"
$ws.Range("B412").Value = "San Nicolás De Los Garza"
$ws.Range("B418").Value = "Guevea De Humboldt"
$ws.Range("B419").Value = "Heroica Ciudad De Ejutla De Crespo"
$ws.Range("B420").Value = "Heroica Ciudad De Juchitán De Zaragoza"
$ws.Range("B421").Value = "Huajuapan De León"
$ws.Range("B422").Value = "Ixtlán De Juárez"
$ws.Range("B425").Value = "Miahuatlán De Porfirio Díaz"
$ws.Range("B426").Value = "Oaxaca De Juárez"
$ws.Range("B429").Value = "San Francisco Del Mar"
$ws.Range("B454").Value = "Tepelmeme Villa De Morelos"
$ws.Range("B455").Value = "Tlacolula De Matamoros"
$ws.Range("B456").Value = "Villa De Tututepec De Melchor Ocampo"
$ws.Range("B465").Value = "Cuapiaxtla De Madero"
$ws.Range("B466").Value = "Cuayuca De Andrade"
$ws.Range("B476").Value = "Palmar De Bravo"
$ws.Range("B484").Value = "San Salvador El Seco"
$ws.Range("B489").Value = "Tepexi De Rodríguez"
$ws.Range("B490").Value = "Tetela De Ocampo"
$ws.Range("B504").Value = "Jalpan De Serra"
$ws.Range("B506").Value = "Pinal De Amoles"
$ws.Range("B509").Value = "San Juan Del Río"
$ws.Range("B522").Value = "Mexquitic De Carmona"
$ws.Range("B524").Value = "San Ciro De Acosta"
$ws.Range("B527").Value = "Santa María Del Río"
$ws.Range("B532").Value = "Villa De Ramos"
$ws.Range("B559").Value = "Nacozari De García"
$ws.Range("B586").Value = "Contla De Juan Cuamatzi"
$ws.Range("B587").Value = "Ixtacuixtla De Mariano Matamoros"
$ws.Range("B589").Value = "Tepetitla De Lardizábal"
$ws.Range("A596").Value = "Veracruz De Ignacio De La Llave"
$ws.Range("B599").Value = "Boca Del Río"
$ws.Range("B601").Value = "Cazones De Herrera"
$ws.Range("B605").Value = "Cosamaloapan De Carpio"
$ws.Range("B616").Value = "Martínez De La Torre"
$ws.Range("B625").Value = "Ozuluama De Mascareñas"
$ws.Range("B628").Value = "Poza Rica De Hidalgo"
$ws.Range("B647").Value = "Cañitas De Felipe Pescador"
$ws.Range("B649").Value = "Concepción Del Oro"
$ws.Range("B665").Value = "Nochistlán De Mejía"
$ws.Range("B666").Value = "Noria De Ángeles"
$ws.Range("B673").Value = "Tlaltenango De Sánchez Román"
$ws.Range("B675").Value = "Villa De Cos"
$ws.Range("A679").Value = "Total"

$ws.Rows("681:685").Delete()

